$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 515, shifting existing rows 515:575 down to 516:576
$ws.Rows.Item(515).Insert()

# Populate the new row 515 with the new record's data.
$ws.Cells.Item(515, 1).Value = 8
$ws.Cells.Item(515, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(515, 3).Value = "Coquimbo"
$ws.Cells.Item(515, 4).Value = 45212
$ws.Cells.Item(515, 5).Value = 4
$ws.Cells.Item(515, 6).Value = 100112021
$ws.Cells.Item(515, 7).Value = "Ají"
$ws.Cells.Item(515, 8).Value = "Inferno"
$ws.Cells.Item(515, 9).Value = "Primera"
$ws.Cells.Item(515, 10).Value = 480
$ws.Cells.Item(515, 11).Value = 28000
$ws.Cells.Item(515, 12).Value = 29000
$ws.Cells.Item(515, 13).Value = 28500
$ws.Cells.Item(515, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(515, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(515, 16).Value = 2850
$ws.Cells.Item(515, 17).Value = 10
$ws.Cells.Item(515, 18).Value = "Hortaliza"
